$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 03:35"

# --- Corea del Sur (row 43): refreshed case counts, same rank ---
$ws.Cells.Item(43,1).Value = "Corea del Sur"
$ws.Cells.Item(43,2).Value = 10936
$ws.Cells.Item(43,3).Value = 27
$ws.Cells.Item(43,4).Value = 9670
$ws.Cells.Item(43,5).Value = 1008
$ws.Cells.Item(43,6).Value = 55
$ws.Cells.Item(43,7).Value = 2
$ws.Cells.Item(43,8).Value = 258

# --- Bolivia's case count rose, moving it above Irak/Grecia/Camerun/Azerbaiyan (rows 70-74) ---
$ws.Cells.Item(70,1).Value = "Bolivia"
$ws.Cells.Item(70,2).Value = 2831
$ws.Cells.Item(70,3).Value = 275
$ws.Cells.Item(70,4).Value = 299
$ws.Cells.Item(70,5).Value = 2410
$ws.Cells.Item(70,6).Value = 3
$ws.Cells.Item(70,7).Value = 4
$ws.Cells.Item(70,8).Value = 122

$ws.Cells.Item(71,1).Value = "Irak"
$ws.Cells.Item(71,2).Value = 2818
$ws.Cells.Item(71,3).Value = 0
$ws.Cells.Item(71,4).Value = 1790
$ws.Cells.Item(71,5).Value = 918
$ws.Cells.Item(71,6).Value = 0
$ws.Cells.Item(71,7).Value = 0
$ws.Cells.Item(71,8).Value = 110

$ws.Cells.Item(72,1).Value = "Grecia"
$ws.Cells.Item(72,2).Value = 2726
$ws.Cells.Item(72,3).Value = 0
$ws.Cells.Item(72,4).Value = 1374
$ws.Cells.Item(72,5).Value = 1201
$ws.Cells.Item(72,6).Value = 32
$ws.Cells.Item(72,7).Value = 0
$ws.Cells.Item(72,8).Value = 151

$ws.Cells.Item(73,1).Value = "Camerun"
$ws.Cells.Item(73,2).Value = 2689
$ws.Cells.Item(73,3).Value = 0
$ws.Cells.Item(73,4).Value = 1524
$ws.Cells.Item(73,5).Value = 1040
$ws.Cells.Item(73,6).Value = 28
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = 125

$ws.Cells.Item(74,1).Value = "Azerbaiyan"
$ws.Cells.Item(74,2).Value = 2589
$ws.Cells.Item(74,3).Value = 0
$ws.Cells.Item(74,4).Value = 1680
$ws.Cells.Item(74,5).Value = 877
$ws.Cells.Item(74,6).Value = 33
$ws.Cells.Item(74,7).Value = 0
$ws.Cells.Item(74,8).Value = 32

# --- Nueva Zelanda (row 88): refreshed active/recovered counts only ---
$ws.Cells.Item(88,4).Value = 1398
$ws.Cells.Item(88,5).Value = 78

# --- Guatemala's case count rose above Somalia (rows 93-94) ---
$ws.Cells.Item(93,1).Value = "Guatemala"
$ws.Cells.Item(93,2).Value = 1114
$ws.Cells.Item(93,3).Value = 62
$ws.Cells.Item(93,4).Value = 111
$ws.Cells.Item(93,5).Value = 977
$ws.Cells.Item(93,6).Value = 5
$ws.Cells.Item(93,7).Value = 0
$ws.Cells.Item(93,8).Value = 26

$ws.Cells.Item(94,1).Value = "Somalia"
$ws.Cells.Item(94,2).Value = 1089
$ws.Cells.Item(94,3).Value = 0
$ws.Cells.Item(94,4).Value = 121
$ws.Cells.Item(94,5).Value = 916
$ws.Cells.Item(94,6).Value = 2
$ws.Cells.Item(94,7).Value = 0
$ws.Cells.Item(94,8).Value = 52

# --- Haiti's case count rose above Santo Tome/Martinica/Islas Feroe/Madagascar (rows 143-147) ---
$ws.Cells.Item(143,1).Value = "Haiti"
$ws.Cells.Item(143,2).Value = 209
$ws.Cells.Item(143,3).Value = 27
$ws.Cells.Item(143,4).Value = 17
$ws.Cells.Item(143,5).Value = 176
$ws.Cells.Item(143,6).Value = 0
$ws.Cells.Item(143,7).Value = 1
$ws.Cells.Item(143,8).Value = 16

$ws.Cells.Item(144,1).Value = "Santo Tome y Principe"
$ws.Cells.Item(144,2).Value = 208
$ws.Cells.Item(144,3).Value = 0
$ws.Cells.Item(144,4).Value = 4
$ws.Cells.Item(144,5).Value = 199
$ws.Cells.Item(144,6).Value = 0
$ws.Cells.Item(144,7).Value = 0
$ws.Cells.Item(144,8).Value = 5

$ws.Cells.Item(145,1).Value = "Martinica"
$ws.Cells.Item(145,2).Value = 187
$ws.Cells.Item(145,3).Value = 0
$ws.Cells.Item(145,4).Value = 83
$ws.Cells.Item(145,5).Value = 90
$ws.Cells.Item(145,6).Value = 4
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 14

$ws.Cells.Item(146,1).Value = "Islas Feroe"
$ws.Cells.Item(146,2).Value = 187
$ws.Cells.Item(146,3).Value = 0
$ws.Cells.Item(146,4).Value = 187
$ws.Cells.Item(146,5).Value = 0
$ws.Cells.Item(146,6).Value = 0
$ws.Cells.Item(146,7).Value = 0
$ws.Cells.Item(146,8).Value = 0

$ws.Cells.Item(147,1).Value = "Madagascar"
$ws.Cells.Item(147,2).Value = 186
$ws.Cells.Item(147,3).Value = 0
$ws.Cells.Item(147,4).Value = 101
$ws.Cells.Item(147,5).Value = 85
$ws.Cells.Item(147,6).Value = 1
$ws.Cells.Item(147,7).Value = 0
$ws.Cells.Item(147,8).Value = 0
